$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) cells whose new value would otherwise be
# auto-converted to a number by Excel, so they stay stored as text strings
# exactly as in the source data (matching the original inlineStr cells).
$textCells = @("D5","D6","D7","D8","D9","D10","D12","D13","D14","D15","D16","D17","D18","D19","D20","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the refreshed cryptos feed.
$ws.Range("D2").Value = "27.126.68"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "1.799.90"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").Value = "1.007"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "308.34"
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("D7").Value = "0.4173"
$ws.Range("E7").Value = "  -1.79%  "
$ws.Range("D8").Value = "0.3550"
$ws.Range("E8").Value = "  -3.21%  "
$ws.Range("D9").Value = "0.07045"
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("D10").Value = "0.8434"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").Value = "2.015.67"
$ws.Range("E11").Value = "  +9.78%  "
$ws.Range("D12").Value = "20.22"
$ws.Range("E12").Value = "  -2.84%  "
$ws.Range("D13").Value = "5.266"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").Value = "6.337"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").Value = "0.06816"
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "79.87"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "0.000008705"
$ws.Range("E18").Value = "  -3.41%  "
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "15.08"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").Value = "27.711.09"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "5.045"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "10.75"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("D24").Value = "2.144.15"
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("D25").Value = "1.951"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "153.10"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").Value = "18.12"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").Value = "5.022"
$ws.Range("E28").Value = "  -4.40%  "
$ws.Range("D29").Value = "112.52"
$ws.Range("E29").Value = "  -2.35%  "
$ws.Range("D30").Value = "1.656"
$ws.Range("E30").Value = "  -10.47%  "
$ws.Range("D31").Value = "0.08866"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "0.7209"
$ws.Range("E32").Value = "  -7.08%  "
$ws.Range("D33").Value = "2.871"
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("D34").Value = "4.338"
$ws.Range("E34").Value = "  -4.92%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Value = "1.078"
$ws.Range("E36").Value = "  -6.38%  "
$ws.Range("D37").Value = "1.077"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("D38").Value = "0.01893"
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("D39").Value = "0.05090"
$ws.Range("E39").Value = "  -5.48%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.1616"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.4933"
$ws.Range("E41").Value = "  -3.74%  "
$ws.Range("D42").Value = "2.593"
$ws.Range("E42").Value = "  -8.20%  "
$ws.Range("D43").Value = "6.196"
$ws.Range("E43").Value = "  -8.40%  "
$ws.Range("D44").Value = "8.075"
$ws.Range("E44").Value = "  -5.28%  "
$ws.Range("D45").Value = "1.007"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").Value = "104.27"
$ws.Range("D47").Value = "10.15"
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("D48").Value = "0.06315"
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("D49").Value = "0.4501"
$ws.Range("E49").Value = "  -4.52%  "
$ws.Range("D50").Value = "1.587"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("D51").Value = "62.49"
$ws.Range("E51").Value = "  -2.91%  "
